# Update countries & provincias Spain
# Applies the daily COVID-data refresh captured by the source diff:
#  - timestamp string in A1 bumped from 17:05 to 17:35
#  - several country rows get refreshed case/death counters
#  - two adjacent-row swaps caused by the re-sort on "Casos totales":
#      Chile leap-frogs Pakistan (rows 22/23)
#      Mozambique leap-frogs Guyana (rows 162/163)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 17:35"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1435110
$ws.Range("C4").Value = 4762
$ws.Range("D4").Value = 310695
$ws.Range("E4").Value = 1039007
$ws.Range("G4").Value = 211
$ws.Range("H4").Value = 85408

# Rows 22/23 - Chile overtakes Pakistan in the ranking
$ws.Range("A22").Value = "Chile"
$ws.Range("B22").Value = 37040
$ws.Range("C22").Value = 2659
$ws.Range("D22").Value = 15655
$ws.Range("E22").Value = 21017
$ws.Range("F22").Value = 555
$ws.Range("G22").Value = 22
$ws.Range("H22").Value = 368

$ws.Range("A23").Value = "Pakistan"
$ws.Range("B23").Value = 35788
$ws.Range("C23").Value = 490
$ws.Range("D23").Value = 9695
$ws.Range("E23").Value = 25323
$ws.Range("F23").Value = 111
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 770

# Row 30 - Singapur
$ws.Range("F30").Value = 20

# Row 74 - Grecia
$ws.Range("B74").Value = 2770
$ws.Range("C74").Value = 10
$ws.Range("E74").Value = 1240
$ws.Range("F74").Value = 24
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 156

# Row 138 - Cabo Verde
$ws.Range("B138").Value = 315
$ws.Range("C138").Value = 26
$ws.Range("D138").Value = 67
$ws.Range("E138").Value = 246

# Rows 162/163 - Mozambique overtakes Guyana in the ranking
$ws.Range("A162").Value = "Mozambique"
$ws.Range("B162").Value = 115
$ws.Range("C162").Value = 11
$ws.Range("D162").Value = 35
$ws.Range("E162").Value = 80
$ws.Range("F162").Value = 0
$ws.Range("H162").Value = 0

$ws.Range("A163").Value = "Guyana"
$ws.Range("B163").Value = 113
$ws.Range("C163").Value = 0
$ws.Range("D163").Value = 41
$ws.Range("E163").Value = 62
$ws.Range("F163").Value = 3
$ws.Range("H163").Value = 10
